# Natmi following Dr Hou advice
# The "Sending cluster" column gains a third value ("sCs"), turning the
# Ntrk3 -> Ptprs sheet into a full cross-product of
# Sending cluster (ECs, FAPs, sCs) x Target cluster (ECs, FAPs, M2, sCs),
# 12 data rows total, with every derived statistic recomputed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A..T per row:
# A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
# E Ligand-expressing cells, F Ligand detection rate,
# G Ligand average expression value, H Ligand total expression value,
# I Ligand derived specificity of average expression value,
# J Ligand derived specificity of total expression value,
# K Receptor-expressing cells, L Receptor detection rate,
# M Receptor average expression value, N Receptor total expression value,
# O Receptor derived specificity of average expression value,
# P Receptor derived specificity of total expression value,
# Q Edge average expression weight, R Edge total expression weight,
# S Edge average expression derived specificity,
# T Edge total expression derived specificity
$data = @(
    @("ECs", "Ntrk3", "Ptprs", "ECs", 1.0, 0.3333333333333333, 0.009008, 0.027024, 0.009775433435787767, 0.009775433435787767, 3.0, 1.0, 4.127188333333333, 12.381565, 0.0561359176022362, 0.05613591760223619, 0.03717771250666666, 0.33459941256, 0.0005487529258775268, 0.0005487529258775267),
    @("ECs", "Ntrk3", "Ptprs", "FAPs", 1.0, 0.3333333333333333, 0.009008, 0.027024, 0.009775433435787767, 0.009775433435787767, 3.0, 1.0, 47.24901333333333, 141.74704, 0.6426570597336346, 0.6426570597336345, 0.4256191121066666, 3.83057200896, 0.006282251309465229, 0.006282251309465227),
    @("ECs", "Ntrk3", "Ptprs", "M2", 1.0, 0.3333333333333333, 0.009008, 0.027024, 0.009775433435787767, 0.009775433435787767, 3.0, 1.0, 8.218848, 24.656544, 0.1117885923419141, 0.1117885923419141, 0.074035382784, 0.666318445056, 0.001092781943318795, 0.001092781943318795),
    @("ECs", "Ntrk3", "Ptprs", "sCs", 1.0, 0.3333333333333333, 0.009008, 0.027024, 0.009775433435787767, 0.009775433435787767, 3.0, 1.0, 13.92629833333334, 41.77889500000001, 0.1894184303222152, 0.1894184303222152, 0.1254480953866667, 1.12903285848, 0.001851647257126218, 0.001851647257126218),
    @("FAPs", "Ntrk3", "Ptprs", "ECs", 3.0, 1.0, 0.07403633333333333, 0.222109, 0.08034383307391152, 0.08034383307391152, 3.0, 1.0, 4.127188333333333, 12.381565, 0.0561359176022362, 0.05613591760223619, 0.305561891176111, 2.750057020584999, 0.004510174793284917, 0.004510174793284916),
    @("FAPs", "Ntrk3", "Ptprs", "FAPs", 3.0, 1.0, 0.07403633333333333, 0.222109, 0.08034383307391152, 0.08034383307391152, 3.0, 1.0, 47.24901333333333, 141.74704, 0.6426570597336346, 0.6426570597336345, 3.498143700817777, 31.48329330736, 0.05163353153100993, 0.05163353153100991),
    @("FAPs", "Ntrk3", "Ptprs", "M2", 3.0, 1.0, 0.07403633333333333, 0.222109, 0.08034383307391152, 0.08034383307391152, 3.0, 1.0, 8.218848, 24.656544, 0.1117885923419141, 0.1117885923419141, 0.608493370144, 5.476440331296, 0.008981524002686288, 0.008981524002686288),
    @("FAPs", "Ntrk3", "Ptprs", "sCs", 3.0, 1.0, 0.07403633333333333, 0.222109, 0.08034383307391152, 0.08034383307391152, 3.0, 1.0, 13.92629833333334, 41.77889500000001, 0.1894184303222152, 0.1894184303222152, 1.031052065506111, 9.279468589555002, 0.0152186027469304, 0.0152186027469304),
    @("sCs", "Ntrk3", "Ptprs", "ECs", 3.0, 1.0, 0.8384493333333333, 2.515348, 0.9098807334903006, 0.9098807334903007, 3.0, 1.0, 4.127188333333333, 12.381565, 0.0561359176022362, 0.05613591760223619, 3.460438306624444, 31.14394475961999, 0.05107698988307375, 0.05107698988307375),
    @("sCs", "Ntrk3", "Ptprs", "FAPs", 3.0, 1.0, 0.8384493333333333, 2.515348, 0.9098807334903006, 0.9098807334903007, 3.0, 1.0, 47.24901333333333, 141.74704, 0.6426570597336346, 0.6426570597336345, 39.6159037299911, 356.54313356992, 0.5847412768931595, 0.5847412768931594),
    @("sCs", "Ntrk3", "Ptprs", "M2", 3.0, 1.0, 0.8384493333333333, 2.515348, 0.9098807334903006, 0.9098807334903007, 3.0, 1.0, 8.218848, 24.656544, 0.1117885923419141, 0.1117885923419141, 6.891087626367999, 62.019788637312, 0.101714286395909, 0.101714286395909),
    @("sCs", "Ntrk3", "Ptprs", "sCs", 3.0, 1.0, 0.8384493333333333, 2.515348, 0.9098807334903006, 0.9098807334903007, 3.0, 1.0, 13.92629833333334, 41.77889500000001, 0.1894184303222152, 0.1894184303222152, 11.67649555338445, 105.08845998046, 0.1723481803181586, 0.1723481803181586)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $i + 2
    for ($j = 0; $j -lt $row.Length; $j++) {
        $c = $j + 1
        $ws.Cells.Item($r, $c).Value = $row[$j]
    }
}
